$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3062.4092
$ws.Range("J64").Value = 3145.4707
$ws.Range("L64").Value = 3145.4707
$ws.Range("N64").Value = -3641.4707
$ws.Range("H67").Value = 3062.4092
$ws.Range("J67").Value = 3145.4707
$ws.Range("L67").Value = 3145.4707
$ws.Range("N67").Value = -4861.4707
$ws.Range("H76").Value = 3244.76
$ws.Range("I76").Value = 3193.762
$ws.Range("K76").Value = 3193.762
$ws.Range("M76").Value = -2878.762
$ws.Range("H79").Value = 3244.76
$ws.Range("I79").Value = 3193.762
$ws.Range("K79").Value = 3193.762
$ws.Range("M79").Value = -2101.762
$ws.Range("H129").Value = 814.0714
$ws.Range("J129").Value = 903.3333
$ws.Range("L129").Value = 2709.9999
$ws.Range("N129").Value = -12709.9999
$ws.Range("H137").Value = 3678.9363
$ws.Range("I137").Value = 3559.48
$ws.Range("J137").Value = 3814.682
$ws.Range("K137").Value = 10678.44
$ws.Range("L137").Value = 11444.046
$ws.Range("M137").Value = -8128.440000000001
$ws.Range("N137").Value = -16544.046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2718.081
$ws.Range("I63").Value = 2564.5
$ws.Range("J63").Value = 2943.3333
$ws.Range("K63").Value = 2564.5
$ws.Range("L63").Value = 2943.3333
$ws.Range("M63").Value = -1878.5
$ws.Range("N63").Value = -4315.3333
$ws.Range("H66").Value = 2718.081
$ws.Range("I66").Value = 2564.5
$ws.Range("J66").Value = 2943.3333
$ws.Range("K66").Value = 12822.5
$ws.Range("L66").Value = 14716.6665
$ws.Range("M66").Value = -9390.5
$ws.Range("N66").Value = -21580.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1685.7894
$ws.Range("I105").Value = 1605
$ws.Range("J105").Value = 2116.6667
$ws.Range("K105").Value = 1605
$ws.Range("L105").Value = 2116.6667
$ws.Range("M105").Value = 142
$ws.Range("N105").Value = -5610.6667
$ws.Range("H134").Value = 1424.7106
$ws.Range("I134").Value = 1195.9584
$ws.Range("J134").Value = 1816.8572
$ws.Range("K134").Value = 3587.8752
$ws.Range("L134").Value = 5450.571599999999
$ws.Range("M134").Value = -1052.8752
$ws.Range("N134").Value = -10520.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3010919.5
$ws.Range("I6").Value = 5017000
$ws.Range("J6").Value = 1799
$ws.Range("K6").Value = 5017000
$ws.Range("L6").Value = 1799
$ws.Range("M6").Value = -5016887
$ws.Range("N6").Value = -2025
$ws.Range("H7").Value = 87.94118
$ws.Range("I7").Value = 46.5
$ws.Range("J7").Value = 110.545456
$ws.Range("K7").Value = 46.5
$ws.Range("L7").Value = 110.545456
$ws.Range("M7").Value = 66.5
$ws.Range("N7").Value = -336.545456
$ws.Range("H58").Value = 2144.12
$ws.Range("I58").Value = 2437.1177
$ws.Range("J58").Value = 1521.5
$ws.Range("K58").Value = 2437.1177
$ws.Range("L58").Value = 1521.5
$ws.Range("M58").Value = -2234.1177
$ws.Range("N58").Value = -1927.5
$ws.Range("H62").Value = 2420
$ws.Range("J62").Value = 2433.3333
$ws.Range("L62").Value = 2433.3333
$ws.Range("N62").Value = -3681.3333
$ws.Range("H65").Value = 2420
$ws.Range("J65").Value = 2433.3333
$ws.Range("L65").Value = 12166.6665
$ws.Range("N65").Value = -18406.6665
$ws.Range("H136").Value = 2144.12
$ws.Range("I136").Value = 2437.1177
$ws.Range("J136").Value = 1521.5
$ws.Range("K136").Value = 7311.353099999999
$ws.Range("L136").Value = 4564.5
$ws.Range("M136").Value = -4761.353099999999
$ws.Range("N136").Value = -9664.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 90910264
$ws.Range("I50").Value = 123.833336
$ws.Range("J50").Value = 200002430
$ws.Range("K50").Value = 371.500008
$ws.Range("L50").Value = 600007290
$ws.Range("M50").Value = 109.499992
$ws.Range("N50").Value = -600008252
$ws.Range("H53").Value = 90910264
$ws.Range("I53").Value = 123.833336
$ws.Range("J53").Value = 200002430
$ws.Range("K53").Value = 371.500008
$ws.Range("L53").Value = 600007290
$ws.Range("M53").Value = 109.499992
$ws.Range("N53").Value = -600008252
$ws.Range("H55").Value = 2999.0908
$ws.Range("J55").Value = 2999.0908
$ws.Range("L55").Value = 8997.2724
$ws.Range("N55").Value = -9351.2724
$ws.Range("H63").Value = 4428.5
$ws.Range("I63").Value = 3233.3333
$ws.Range("J63").Value = 8014
$ws.Range("K63").Value = 9699.999899999999
$ws.Range("L63").Value = 24042
$ws.Range("M63").Value = -8950.999899999999
$ws.Range("N63").Value = -25540
$ws.Range("H64").Value = 2414.4
$ws.Range("I64").Value = 2012
$ws.Range("J64").Value = 2476.3076
$ws.Range("K64").Value = 6036
$ws.Range("L64").Value = 7428.9228
$ws.Range("M64").Value = -5766
$ws.Range("N64").Value = -7968.9228
$ws.Range("H66").Value = 4428.5
$ws.Range("I66").Value = 3233.3333
$ws.Range("J66").Value = 8014
$ws.Range("K66").Value = 29099.9997
$ws.Range("L66").Value = 72126
$ws.Range("M66").Value = -25355.9997
$ws.Range("N66").Value = -79614
$ws.Range("H67").Value = 2414.4
$ws.Range("I67").Value = 2012
$ws.Range("J67").Value = 2476.3076
$ws.Range("K67").Value = 6036
$ws.Range("L67").Value = 7428.9228
$ws.Range("M67").Value = -5100
$ws.Range("N67").Value = -9300.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 60000
$ws.Range("J20").Value = 60000
$ws.Range("L20").Value = 60000
$ws.Range("N20").Value = -60490
$ws.Range("H70").Value = 4134.826
$ws.Range("I70").Value = 4131.6
$ws.Range("J70").Value = 4140.875
$ws.Range("K70").Value = 4131.6
$ws.Range("L70").Value = 4140.875
$ws.Range("M70").Value = -3861.6
$ws.Range("N70").Value = -4680.875
$ws.Range("H73").Value = 4134.826
$ws.Range("I73").Value = 4131.6
$ws.Range("J73").Value = 4140.875
$ws.Range("K73").Value = 4131.6
$ws.Range("L73").Value = 4140.875
$ws.Range("M73").Value = -3195.6
$ws.Range("N73").Value = -6012.875
$ws.Range("H80").Value = 2737
$ws.Range("J80").Value = 2618.4
$ws.Range("L80").Value = 2618.4
$ws.Range("N80").Value = -4614.4
$ws.Range("H83").Value = 2737
$ws.Range("J83").Value = 2618.4
$ws.Range("L83").Value = 13092
$ws.Range("N83").Value = -23076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 459.04544
$ws.Range("I55").Value = 409.1875
$ws.Range("J55").Value = 592
$ws.Range("K55").Value = 409.1875
$ws.Range("L55").Value = 592
$ws.Range("M55").Value = -236.1875
$ws.Range("N55").Value = -938
$ws.Range("H136").Value = 3024.7556
$ws.Range("I136").Value = 2351.3076
$ws.Range("K136").Value = 7053.9228
$ws.Range("M136").Value = -4503.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1252101.8
$ws.Range("I122").Value = 9999999
$ws.Range("J122").Value = 2402.1428
$ws.Range("K122").Value = 29999997
$ws.Range("L122").Value = 7206.428400000001
$ws.Range("M122").Value = -29997547
$ws.Range("N122").Value = -12106.4284
$ws.Range("H126").Value = 2002161.8
$ws.Range("I126").Value = 5001249.5
$ws.Range("J126").Value = 2770
$ws.Range("K126").Value = 15003748.5
$ws.Range("L126").Value = 8310
$ws.Range("M126").Value = -15001278.5
$ws.Range("N126").Value = -13250
